# Kent's feedback: reposition/resize the "F" callout label on the cross-section
# diagram (slide 2) and swap its Swell-Braille glyph text for a plain "B".
#
# Note: Shape.Left/Top/Width/Height are expressed in points (1 pt = 12700 EMU).
# The literals below are EMU/12700 (nudged a hair above the exact quotient so
# the host's single-precision round-trip lands back on the exact target EMU
# instead of one EMU short):
#   off x=8307571 -> 8456426 EMU   ext cx=649537 -> 407484 EMU
#   off y=3044455 -> 3129516 EMU   ext cy=461665 -> 461665 EMU (unchanged)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item("TextBox 47")

$shp.Left = 665.8603219606299
$shp.Top = 246.41858767716533
$shp.Width = 32.08535633070866
$shp.Height = 36.351576803149605

$shp.TextFrame.TextRange.Text = "B"
